# TrialsSetup.xlsx update - 2026-02-09 12:00
# Sheet1 "Days remaining" figures refreshed for two trials:
#   - REJOICE (MK-5909-003)      row 7 -> 10 (was 12)
#   - REMASTER (CLOU)            row 9 -> 30 (was 32)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 10
$ws.Range("B9").Value = 30
